# "New plan from Director" -- add sign-off initials to several section
# leads, drop the stray _GoBack bookmark from the Sound Lead Q&A, and
# insert the Director's new step-by-step plan (with a numbered list)
# right after the "Director will use this information..." paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: append `$text` to the end of `$para` as a NEW run rather than
# merging into the paragraph's trailing run. Word's Range.InsertAfter
# happily folds identical-formatting text into the adjacent run when a
# paragraph has only a single run, which does not match what a real
# edit (typed at a different time / by a different author) looks like
# in the underlying XML. Splitting the paragraph and then deleting the
# paragraph mark keeps the two pieces of text as distinct <w:r> runs.
# ---------------------------------------------------------------------
function Append-NewRun($para, [string]$text) {
    $para.Range.InsertParagraphAfter()
    $nextPara = $para.Next()
    $nextPara.Range.Text = $text
    $mark = $d.Range($para.Range.End - 1, $para.Range.End)
    $mark.Delete()
}

# 1) Design Lead paragraph gets " -Noah"
$rng = $d.Content
$rng.Find.Execute("Document level design and player behavior", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" -Noah")

# 2) Sound Lead paragraph gets " -Nick"
$rng = $d.Content
$rng.Find.Execute("Document all audio elements in the project.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" -Nick")

# 3) Drop the stray _GoBack bookmark that currently sits at the end of
#    the "They start on awake and end on death..." paragraph; it is
#    re-created later, further down, once the new content is in place.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 4) Tech Lead paragraph gets " -Noah"
$rng = $d.Content
$rng.Find.Execute("Document systems and interactions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" -Noah")

# 5) Director paragraph gets " –Zeke" (en dash) as a separate run
$rng = $d.Content
$rng.Find.Execute("Director will use this information to decide a plan of action.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$directorPara = $rng.Paragraphs(1)
Append-NewRun $directorPara " –Zeke"

# 6) Insert the new "Step" paragraph + the 5 numbered steps right after
#    the Director paragraph (before the Producer paragraph).
$insertAfter = $directorPara

$insertAfter.Range.InsertParagraphAfter()
$stepPara = $insertAfter.Next()
$stepPara.Range.Text = " Step"

$stepPara.Range.InsertParagraphAfter()
$item1 = $stepPara.Next()
$item1.Range.Text = "Each get our own copy of the broken game and start finding what is the real assest for the game"
$item1.Range.ListFormat.ApplyNumberDefault()

$item1.Range.InsertParagraphAfter()
$item2 = $item1.Next()
$item2.Range.Text = "Make a local broken copy of the game for github"

$item2.Range.InsertParagraphAfter()
$item3 = $item2.Next()
$item3.Range.Text = "Sort out the real stuff from the broken stuff into different folders"

$item3.Range.InsertParagraphAfter()
$item4 = $item3.Next()
$item4.Range.Text = "Individually in personal sandboxes’’ start putting game together (For sound and others that require the others, they can get a fixed or pieced items to test on)"

$item4.Range.InsertParagraphAfter()
$item5 = $item4.Next()
$item5.Range.Text = "After done we put all pieces in final copy (connect sprites to prefabs first) so there’s everything in the game then start attaching the physics to the objects and then add sound once things are in."

# 7) Re-create the _GoBack bookmark right after "...then start attaching"
#    (before " the physics to the objects...") inside the last step item.
$rng = $d.Content
$rng.Find.Execute("then start attaching", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)

# 8) Producer paragraph gets " -Nick"
$rng = $d.Content
$rng.Find.Execute("Producer will consolidate this information to be used by the team.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$producerPara = $rng.Paragraphs(1)
Append-NewRun $producerPara " -Nick"
